$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 9 new rows of accelerometer data right after the header,
#     pushing the existing data down. ---
#
# Inserting directly below the header row causes Excel to inherit the
# header's bold/centered formatting for the new rows, so instead we insert
# the blank rows further down (row 11, inside the homogeneous data block)
# where neighbouring rows are unstyled, then cut/paste the original data
# down into the freshly inserted rows. This keeps formatting identical to
# the rest of the plain data rows.
$ws.Rows.Item(11).Resize(9).Insert() | Out-Null
$ws.Range("A2:C10").Cut($ws.Range("A11:C19")) | Out-Null

$newTopRows = @(
    @(-0.5562429428100586, 1.515548229217529, 0.1288182139396667),
    @(-0.7228193283081055, 1.569920063018799, 0.3090478777885437),
    @(-0.7100648880004883, 1.677208662033081, 0.4076560139656067),
    @(-0.8624534606933594, 1.662384271621704, 0.1703254878520965),
    @(-0.2973442077636719, 1.564767122268677, 0.0300358235836029),
    @(-0.7404184341430664, 1.651389360427856, 0.2246546447277069),
    @(-0.7922754287719727, 1.612479209899902, 0.1184005141258239),
    @(-0.5356760025024414, 1.573039531707764, 0.1370119750499725),
    @(-0.6788949966430664, 1.517318725585938, 0.1871603727340698)
)

$r = 2
foreach ($row in $newTopRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# --- Append one new row of accelerometer data at the end. ---
$ws.Cells.Item(31, 1).Value = -0.1529102325439453
$ws.Cells.Item(31, 2).Value = 1.506775379180908
$ws.Cells.Item(31, 3).Value = 0.2578078508377075
